$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The paragraph right after the "Ancillary Structures" heading (styled
# "FirstParagraph") currently holds nothing but a tiny inline placeholder
# picture. Swap that picture out for a hyperlink run whose visible text is
# the image's own source URL (the .jpg it used to point at).
# ---------------------------------------------------------------------------
$jpgUrl = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/CCI03_Ancillary_Structures.jpg?h=100%25&w=100%25"

$shape = $d.InlineShapes.Item(1)
$pictureStart = $shape.Range.Start

# Remove the picture, then insert the hyperlink exactly where it used to sit.
$shape.Delete()
$insertAt = $d.Range($pictureStart, $pictureStart)
$d.Hyperlinks.Add($insertAt, $jpgUrl, "", "", $jpgUrl)

Write-Output "Replaced ancillary-structures picture with hyperlink to $jpgUrl"
